# Apply the "household indicators" edit to the MOH 515 (Post Outbreak) survey
# workbook: re-case a handful of existing question labels, rename the first
# group's internal name, and append a brand-new "Household Indicators" group
# (5 integer questions) just before the form's closing "end group" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$xlPasteFormats = -4122

# --- Re-case a handful of existing question labels / the group name ---
$ws.Range("B2").Value = "form_summary"

$ws.Range("C3").Value = "What Is Your Name?"
$ws.Range("C4").Value = "What Is Your Area?"
$ws.Range("C5").Value = "What Is Your Linked Facility?"
$ws.Range("C6").Value = "What County Do You Belong To?"
$ws.Range("C7").Value = "How Many CHPs Are In Your Area?"
$ws.Range("C8").Value = "How Many CHPs Submitted Monthly Report?"

# --- Append the new "Household Indicators" group in rows 10-16, directly
#     below the existing "end group" row (row 9), which stays put. Cells are
#     written straight to their final row (no Rows.Insert(), which would
#     smear row 8/9's per-cell formatting into the untouched columns D/F/G/H
#     of the new rows). Formatting is carried over with a Copy/PasteSpecial
#     of just the (already-used) cell styles, then the value is set.

function Set-StyledCell($targetAddr, $styleSourceAddr, $value) {
    $ws.Range($styleSourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($targetAddr).Value = $value
}

# Row 10: begin group / household_indicators / Household Indicators
Set-StyledCell "A10" "A9" "begin group"
Set-StyledCell "B10" "B8" "household_indicators"
Set-StyledCell "C10" "C8" "Household Indicators"

# Row 11: integer / total_households / Total Households In The Area? / yes / numbers
Set-StyledCell "A11" "A8" "integer"
Set-StyledCell "B11" "B8" "total_households"
Set-StyledCell "C11" "C8" "Total Households In The Area?"
Set-StyledCell "D11" "D8" "yes"
Set-StyledCell "F11" "F8" "numbers"

# Row 12: integer / new_households / Number Of New Households Registered This Month? / yes / numbers
Set-StyledCell "A12" "A8" "integer"
Set-StyledCell "B12" "B8" "new_households"
Set-StyledCell "C12" "C8" "Number Of New Households Registered This Month?"
Set-StyledCell "D12" "D8" "yes"
Set-StyledCell "F12" "F8" "numbers"

# Row 13: integer / new_households_visited / Number Of New Households Visited This Month? / yes / numbers
Set-StyledCell "A13" "A8" "integer"
Set-StyledCell "B13" "B8" "new_households_visited"
Set-StyledCell "C13" "C8" "Number Of New Households Visited This Month?"
Set-StyledCell "D13" "D8" "yes"
Set-StyledCell "F13" "F8" "numbers"

# Row 14: integer / new_households_with_clean_water / ... / yes  (no appearance)
Set-StyledCell "A14" "A8" "integer"
Set-StyledCell "B14" "B8" "new_households_with_clean_water"
Set-StyledCell "C14" "C8" "Number Of New Households Visited This Month With Clean Water Access?"
Set-StyledCell "D14" "D8" "yes"

# Row 15: integer / new_households_with_latrines / ... / yes (label cell keeps the shaded style)
Set-StyledCell "A15" "A8" "integer"
Set-StyledCell "B15" "B8" "new_households_with_latrines"
Set-StyledCell "C15" "H8" "Number Of New Households Visited This Month With Latrines/Toilets?"
Set-StyledCell "D15" "D8" "yes"

# Row 16: closing end group for the new section
Set-StyledCell "A16" "A9" "end group"

# --- Column widths: a new, narrower column B, and a much wider column C to
#     fit the longer household-indicator question labels. ---
$ws.Columns("B").ColumnWidth = 17.92
$ws.Columns("C").ColumnWidth = 59.25
